$wb = $excel.ActiveWorkbook
$sEV = $wb.Worksheets.Item("Earned Value")
$sHI = $wb.Worksheets.Item("Horas insumidas")

# ---------------------------------------------------------------------------
# 1) "Horas insumidas": append the two new time-tracking rows for the new
#    story S-01030 (rows 136/137), matching the date-number-format of the
#    existing block (row 135).
# ---------------------------------------------------------------------------
$sHI.Range("B136").Value = 40501
$sHI.Range("C136").Value = "Sergio"
$sHI.Range("D136").Value = "Se corrige el bug de prioridad media de validación de rangos para las métricas"
$sHI.Range("E136").Value = "S-01030"
$sHI.Range("F136").Value = 3
$sHI.Range("B136").NumberFormat = $sHI.Range("B135").NumberFormat()

$sHI.Range("B137").Value = 40502
$sHI.Range("C137").Value = "Sergio"
$sHI.Range("D137").Value = "Se corrige el bug de prioridad media de validación de rangos para las métricas"
$sHI.Range("E137").Value = "S-01030"
$sHI.Range("F137").Value = 6
$sHI.Range("B137").NumberFormat = $sHI.Range("B135").NumberFormat()

# ---------------------------------------------------------------------------
# 2) "Earned Value": insert a new row for story S-01030 right after the
#    existing "Métricas para agentes: 3 métricas" row (old row 22), pushing
#    the following rows down by one. The inherited row formatting already
#    matches (columns A/B styled, E styled as the numeric PV column).
# ---------------------------------------------------------------------------
$sEV.Rows.Item(23).Insert()

$sEV.Range("A23").Value = "S-01030"
$sEV.Range("B23").Value = "Arreglar todos los bugs de prioridad media/alta que figuran en el informe de avance"
$sEV.Range("C23").Value = "Completada"
$sEV.Range("D23").Value = 100
$sEV.Range("E23").Value = 20
$sEV.Range("F23").Value = 20
$sEV.Range("G23").Formula = "=SUMIF('Horas insumidas'!`$E`$6:`$E`$150,A23,'Horas insumidas'!`$F`$6:`$F`$150)"
$sEV.Range("H23").Formula = "=F23-G23"
$sEV.Range("I23").Formula = "=F23-E23"
$sEV.Range("J23").Formula = "=F23/E23"
$sEV.Range("K23").Formula = "=F23/G23"

# Fix the id that used to be inconsistently typed ("S01024" -> "S-01024").
$sEV.Range("A22").Value = "S-01024"

# Add the tracking link on the new story's description cell.
$sEV.Hyperlinks.Add($sEV.Range("B23"), "https://www1.v1host.com/Team152/assetdetail.v1?oid=Story%3a1191", "", "", "https://www1.v1host.com/Team152/assetdetail.v1?oid=Story%3a1191") | Out-Null

# Widen every SUMIF lookup range on sheet "Earned Value" (rows 2-22) so the
# "Horas insumidas" table's extra rows (up to 150) are included.
for ($r = 2; $r -le 22; $r++) {
    $sEV.Range("G$r").Formula = "=SUMIF('Horas insumidas'!`$E`$6:`$E`$150,A$r,'Horas insumidas'!`$F`$6:`$F`$150)"
}

# ---------------------------------------------------------------------------
# 3) View state: the active tab moves from "Horas insumidas" to
#    "Earned Value", selections move, and "Horas insumidas" scrolls a bit.
# ---------------------------------------------------------------------------
$sHI.Range("B138").Select() | Out-Null
$sEV.Range("K23").Select() | Out-Null
